$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = ""

$ws.Range("C3").Value = ""
$ws.Range("C4").Value = ""
$ws.Range("C5").Value = ""

$ws.Range("C6").Value = "8.300,01 TL - 199,41 TL"
$ws.Range("J6").Value = "8.300,01 TL - 199,41 TL"
$ws.Range("K6").Value = "7,97 TL - 15,96 TL - 199,41 TL"

$ws.Range("D7").Value = ""

$ws.Range("C8").Value = ""
$ws.Range("C9").Value = ""
$ws.Range("C10").Value = ""
$ws.Range("C11").Value = ""

$ws.Range("K12").Value = "WU: ,USD–; Diğer: 529 TL–4.454,74 TL"

$ws.Range("C13").Value = "Hesaba: Asgari 0 TL | Azami 0,94 TL"
$ws.Range("J13").Value = "Hesaba: Asgari 1 TL | Azami 995,5 TL"
$ws.Range("K13").Value = "Hesaba: Asgari 1 TL | Azami 865,75 TL"

$ws.Range("C14").Value = "40.000 TL - 2.485,72 TL"
$ws.Range("J14").Value = "1.554,97 TL - 7.784 TL"
$ws.Range("K14").Value = "1.196,51 TL - 5.583,74 TL"

$ws.Range("C15").Value = " Asgari Tutar: 390,48 TL Azami Tutar: 390,48 TL"

$ws.Range("C17").Value = " Asgari Tutar: 257,15 TL Azami Tutar: 257,15 TL"

$ws.Range("C20").Value = "114,29 TL"

$ws.Range("C21").Value = "%0,57 Asgari Tutar: 590,48 TL Azami Tutar: 590,48 TL / 3.295,24 TL"

$ws.Range("C23").Value = "85,72 TL"

$ws.Range("C24").Value = "600 TL"

$ws.Range("C25").Value = "495,24 TL"
